# #5: property boat&car done
# Rework the "汽車" (car) sheet: row 1 was (incorrectly) holding data-looking
# values instead of column headers, and row 2 only carried 7 of the columns
# that every other property sheet has. Fix the header row and extend both
# rows out to the full common column set (property_category .. index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

function Set-HeaderCell($col, $text) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
}

# --- Row 1: proper header labels (B1:N1) ---
Set-HeaderCell 2  "name"
Set-HeaderCell 3  "capacity"
Set-HeaderCell 4  "owner"
Set-HeaderCell 5  "register_date"
Set-HeaderCell 6  "register_reason"
Set-HeaderCell 7  "acquire_value"
Set-HeaderCell 8  "property_category"
Set-HeaderCell 9  "category"
Set-HeaderCell 10 "date"
Set-HeaderCell 11 "legislator_name"
Set-HeaderCell 12 "legislator_id"
Set-HeaderCell 13 "source_file"
Set-HeaderCell 14 "index"

# --- Row 2: extend the data row with the shared trailing columns ---
$ws.Cells.Item(2, 8).Value  = "land"
$ws.Cells.Item(2, 9).Value  = "normal"
# force text so "2011-11-28" isn't auto-converted into a date serial number
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2011-11-28"
$ws.Cells.Item(2, 11).Value = "林德福"
$ws.Cells.Item(2, 12).Value = 908
$ws.Cells.Item(2, 13).Value = "tmp2e4a1"
$ws.Cells.Item(2, 14).Value = 40
